$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at the top. This shifts the existing header row
#    (Field/Data Type/Source/Description) from row 1 -> row 2, and the
#    raw-field data rows from rows 2-8 -> rows 3-9.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).Insert()

# New section title in the now-empty row 1.
$ws.Range("A1").Value = "Raw Data Fields"

# ---------------------------------------------------------------------------
# 2. Update the raw "open/high/low/close" rows (now rows 5-8) from INT to
#    FLOAT and reword their descriptions.
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "FLOAT"
$ws.Range("D5").Value = "Opening Price in USD"

$ws.Range("B6").Value = "FLOAT"
$ws.Range("D6").Value = "Highest price during trading day in USD"

$ws.Range("B7").Value = "FLOAT"
$ws.Range("D7").Value = "Lowest price during trading day in USD"

$ws.Range("B8").Value = "FLOAT"
$ws.Range("D8").Value = "Closing Price in USD"

# ---------------------------------------------------------------------------
# 3. New "Processed Data Fields" section (rows 11-16).
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "Processed Data Fields"

$ws.Range("A12").Value = "Field"
$ws.Range("B12").Value = "Data Type"
$ws.Range("C12").Value = "Source"
$ws.Range("D12").Value = "Description"

$ws.Range("A13").Value = "daily_return"
$ws.Range("B13").Value = "FLOAT"
$ws.Range("C13").Value = "data/staged/processed_stocks.csv"
$ws.Range("D13").Value = "Log daily return: ln(close_t/close_t-1)"

$ws.Range("A14").Value = "volatility_14d"
$ws.Range("B14").Value = "FLOAT"
$ws.Range("C14").Value = "data/staged/processed_stocks.csv"
$ws.Range("D14").Value = "14-day rolling volatility (annualized using √252)"

$ws.Range("A15").Value = "simple_return"
$ws.Range("B15").Value = "FLOAT"
$ws.Range("C15").Value = "data/staged/processed_stocks.csv"
$ws.Range("D15").Value = "Simple daily return: (close_t/close_t-1) - 1"

$ws.Range("A16").Value = "cumulative_return"
$ws.Range("B16").Value = "FLOAT"
$ws.Range("C16").Value = "data/staged/processed_stocks.csv"
$ws.Range("D16").Value = "Cumulative return from start of period"

# ---------------------------------------------------------------------------
# 4. New "DATA QUALITY NOTES" section (rows 22-26).
# ---------------------------------------------------------------------------
$ws.Range("A22").Value = "DATA QUALITY NOTES"
$ws.Range("A23").Value = "- Missing Values Expected: First row per ticker has NaN daily_return (no previous price)"
$ws.Range("A24").Value = "- Rolling Window: First 13 rows per ticker have NaN volatility_14d (insufficient data for 14-day window)"
$ws.Range("A25").Value = "- Data Range: Approximately 250 trading days per ticker (one year of market data)"
$ws.Range("A26").Value = "- Annualization Factor: Volatility multiplied by √252 to convert daily to annual measure"

# ---------------------------------------------------------------------------
# 5. Formatting: section-title rows get a bold, 14pt font + taller row, the
#    two new "Field/Data Type/Source/Description" header rows match the
#    original (bold, default size) header styling.
# ---------------------------------------------------------------------------
foreach ($r in @(1, 11, 22)) {
    $cell = $ws.Range("A$r")
    $cell.Font.Bold = $true
    $cell.Font.Size = 14
    $ws.Rows.Item($r).RowHeight = 18.5
}

$ws.Range("A12:D12").Font.Bold = $true

# ---------------------------------------------------------------------------
# 6. Column widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 25.666666666666668
$ws.Columns.Item(3).ColumnWidth = 31.666666666666668
$ws.Columns.Item(4).ColumnWidth = 43.333333333333336

# ---------------------------------------------------------------------------
# 7. Selection.
# ---------------------------------------------------------------------------
$ws.Range("A26").Select()
